$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Intakes")

# Remove the "intake_funding_source" column (column J) from the Intakes sheet.
# Deleting the entire column shifts subsequent columns (K..O) left by one.
$ws.Range("J:J").EntireColumn.Delete()
